# Update cryptos list data (Price and Volume(1h) columns) for rows 2-51
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column values are stored as text (matches original inline-string cells),
# since several values (e.g. "1.000", "13.10", "244.06") would otherwise be
# auto-coerced to numbers and lose their exact textual formatting.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '30.346.68'
$ws.Range("E2").Value = '  -1.25%  '
$ws.Range("D3").Value = '1.869.32'
$ws.Range("E3").Value = '  -1.12%  '
$ws.Range("D4").Value = '0.9997'
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '244.06'
$ws.Range("E5").Value = '  -2.22%  '
$ws.Range("E6").Value = '  +0.02%  '
$ws.Range("D7").Value = '0.4720'
$ws.Range("E7").Value = '  -0.78%  '
$ws.Range("D8").Value = '0.2876'
$ws.Range("E8").Value = '  -2.23%  '
$ws.Range("D9").Value = '0.06485'
$ws.Range("E9").Value = '  -0.85%  '
$ws.Range("D10").Value = '21.82'
$ws.Range("E10").Value = '  -1.19%  '
$ws.Range("D11").Value = '100.02'
$ws.Range("E11").Value = '  +3.21%  '
$ws.Range("D12").Value = '0.07795'
$ws.Range("E12").Value = '  +0.66%  '
$ws.Range("D13").Value = '1.871.96'
$ws.Range("E13").Value = '  -1.00%  '
$ws.Range("D14").Value = '0.7270'
$ws.Range("E14").Value = '  -2.03%  '
$ws.Range("D15").Value = '5.167'
$ws.Range("E15").Value = '  -1.70%  '
$ws.Range("D16").Value = '283.23'
$ws.Range("E16").Value = '  +2.53%  '
$ws.Range("D17").Value = '30.325.39'
$ws.Range("E17").Value = '  -1.30%  '
$ws.Range("D18").Value = '13.10'
$ws.Range("E18").Value = '  -0.99%  '
$ws.Range("D19").Value = '1.000'
$ws.Range("E19").Value = '  +0.05%  '
$ws.Range("D20").Value = '0.000007482'
$ws.Range("E20").Value = '  -1.24%  '
$ws.Range("D21").Value = '2.113.84'
$ws.Range("E21").Value = '  -0.95%  '
$ws.Range("E22").Value = '  +0.04%  '
$ws.Range("D23").Value = '5.276'
$ws.Range("E23").Value = '  -1.13%  '
$ws.Range("D24").Value = '6.287'
$ws.Range("E24").Value = '  +0.57%  '
$ws.Range("D25").Value = '163.10'
$ws.Range("E25").Value = '  -0.65%  '
$ws.Range("D26").Value = '9.040'
$ws.Range("E26").Value = '  -2.37%  '
$ws.Range("D27").Value = '19.01'
$ws.Range("E27").Value = '  +0.82%  '
$ws.Range("D28").Value = '1.891'
$ws.Range("E28").Value = '  -1.98%  '
$ws.Range("D29").Value = '0.09646'
$ws.Range("E29").Value = '  -0.89%  '
$ws.Range("D30").Value = '1.319'
$ws.Range("E30").Value = '  -1.97%  '
$ws.Range("D31").Value = '1.482'
$ws.Range("E31").Value = '  -1.69%  '
$ws.Range("D32").Value = '4.229'
$ws.Range("E32").Value = '  -1.90%  '
$ws.Range("D33").Value = '4.140'
$ws.Range("E33").Value = '  -1.19%  '
$ws.Range("D34").Value = '0.04805'
$ws.Range("E34").Value = '  -1.75%  '
$ws.Range("D35").Value = '1.124'
$ws.Range("E35").Value = '  -0.34%  '
$ws.Range("D36").Value = '0.6883'
$ws.Range("E36").Value = '  -1.89%  '
$ws.Range("D37").Value = '2.716'
$ws.Range("E37").Value = '  -0.28%  '
$ws.Range("D38").Value = '0.01896'
$ws.Range("E38").Value = '  -1.20%  '
$ws.Range("E39").Value = '  +1.53%  '
$ws.Range("D40").Value = '75.71'
$ws.Range("E40").Value = '  -0.32%  '
$ws.Range("D41").Value = '6.269'
$ws.Range("E41").Value = '  -1.42%  '
$ws.Range("D42").Value = '1.954'
$ws.Range("E42").Value = '  -3.94%  '
$ws.Range("D43").Value = '0.4220'
$ws.Range("E43").Value = '  -0.99%  '
$ws.Range("D44").Value = '0.9989'
$ws.Range("E44").Value = '  -0.09%  '
$ws.Range("D45").Value = '0.8232'
$ws.Range("E45").Value = '  -2.43%  '
$ws.Range("D46").Value = '100.73'
$ws.Range("E46").Value = '  -1.66%  '
$ws.Range("D47").Value = '9.794'
$ws.Range("E47").Value = '  +3.95%  '
$ws.Range("D48").Value = '7.004'
$ws.Range("E48").Value = '  -1.08%  '
$ws.Range("D49").Value = '35.01'
$ws.Range("E49").Value = '  -2.00%  '
$ws.Range("D50").Value = '0.05762'
$ws.Range("E50").Value = '  -0.30%  '
$ws.Range("D51").Value = '885.60'
$ws.Range("E51").Value = '  -3.70%  '
